# Fruta / hortaliza, semanal
#
# A new weekly price-report row is inserted above the current row 444 of
# "Sheet1" (Vega Modelo de Temuco - Perejil). Inserting the row pushes the
# existing rows 444-503 down to 445-504 (so the previously-last row, with
# date 45072, ends up at row 504), and the freshly-opened row 444 is filled
# in with this week's reading.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 444:503 down to 445:504, leaving a blank row 444 in place.
$ws.Rows("444:444").Insert()

# Populate the newly inserted row 444 with the new weekly record.
$ws.Range("A444").Value = 10
$ws.Range("B444").Value = 'Vega Modelo de Temuco'
$ws.Range("C444").Value = 'La Araucanía'
$ws.Range("D444").Value2 = 45124
$ws.Range("E444").Value = 9
$ws.Range("F444").Value = 100112044
$ws.Range("G444").Value = 'Perejil'
$ws.Range("H444").Value = 'Sin especificar'
$ws.Range("I444").Value = 'Primera'
$ws.Range("J444").Value = 50
$ws.Range("K444").Value = 4000
$ws.Range("L444").Value = 4000
$ws.Range("M444").Value = 4000
$ws.Range("N444").Value = '$/docena de atados (3 kilos)'
$ws.Range("O444").Value = 'Provincia de Cautín'
$ws.Range("P444").Value = 1333
$ws.Range("Q444").Value = 3
$ws.Range("R444").Value = 'Hortaliza'
